$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = 1145
$ws.Range("D8").Value = 189
$ws.Range("E8").Value = 956
$ws.Range("F8").Value = 7.752255947497949
$ws.Range("G8").Value = 83.49344978165939
$ws.Range("H8").Value = 16.50655021834061
